$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.165316581726074
$ws.Range("B1").Value = 5.136908054351807
$ws.Range("C1").Value = 4.227310180664062
$ws.Range("D1").Value = 5.021109104156494
$ws.Range("E1").Value = 4.630776405334473
